# Apply updated "F" column (人气/热度 count) values as published by the
# regenerated gh-pages output (commit 456a3b4).
#
# Each entry: worksheet name, cell address, expected-old value, new value.
# The old value is only used as a sanity check so we never touch the wrong
# cell if the source workbook differs from what we expect.

$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "展览";     Cell = "F3";  Old = 1752; New = 1753 },
    @{ Sheet = "展览";     Cell = "F5";  Old = 493;  New = 498 },
    @{ Sheet = "展览";     Cell = "F8";  Old = 1253; New = 1256 },
    @{ Sheet = "展览";     Cell = "F9";  Old = 366;  New = 367 },
    @{ Sheet = "展览";     Cell = "F11"; Old = 895;  New = 897 },
    @{ Sheet = "展览";     Cell = "F13"; Old = 194;  New = 195 },
    @{ Sheet = "展览";     Cell = "F14"; Old = 525;  New = 526 },
    @{ Sheet = "展览";     Cell = "F18"; Old = 2989; New = 2991 },
    @{ Sheet = "展览";     Cell = "F19"; Old = 2637; New = 2638 },
    @{ Sheet = "展览";     Cell = "F26"; Old = 5363; New = 5366 },
    @{ Sheet = "展览";     Cell = "F31"; Old = 339;  New = 340 },

    @{ Sheet = "演出";     Cell = "F4";  Old = 1151; New = 1152 },
    @{ Sheet = "演出";     Cell = "F14"; Old = 618;  New = 619 },
    @{ Sheet = "演出";     Cell = "F20"; Old = 617;  New = 618 },
    @{ Sheet = "演出";     Cell = "F25"; Old = 283;  New = 284 },
    @{ Sheet = "演出";     Cell = "F26"; Old = 3980; New = 3982 },
    @{ Sheet = "演出";     Cell = "F34"; Old = 35;   New = 36 },

    @{ Sheet = "本地生活"; Cell = "F5";  Old = 2502; New = 2504 },
    @{ Sheet = "本地生活"; Cell = "F6";  Old = 1065; New = 1067 },
    @{ Sheet = "本地生活"; Cell = "F9";  Old = 1356; New = 1357 },
    @{ Sheet = "本地生活"; Cell = "F10"; Old = 372;  New = 374 },

    @{ Sheet = "全部类型"; Cell = "F5";  Old = 2502; New = 2504 },
    @{ Sheet = "全部类型"; Cell = "F6";  Old = 1752; New = 1753 },
    @{ Sheet = "全部类型"; Cell = "F7";  Old = 1065; New = 1067 },
    @{ Sheet = "全部类型"; Cell = "F8";  Old = 1356; New = 1357 },
    @{ Sheet = "全部类型"; Cell = "F9";  Old = 372;  New = 374 },
    @{ Sheet = "全部类型"; Cell = "F12"; Old = 494;  New = 498 },
    @{ Sheet = "全部类型"; Cell = "F15"; Old = 1253; New = 1256 },
    @{ Sheet = "全部类型"; Cell = "F16"; Old = 366;  New = 367 },
    @{ Sheet = "全部类型"; Cell = "F17"; Old = 895;  New = 897 },
    @{ Sheet = "全部类型"; Cell = "F19"; Old = 1151; New = 1152 },
    @{ Sheet = "全部类型"; Cell = "F20"; Old = 1151; New = 1152 },
    @{ Sheet = "全部类型"; Cell = "F21"; Old = 194;  New = 195 },
    @{ Sheet = "全部类型"; Cell = "F22"; Old = 525;  New = 526 },
    @{ Sheet = "全部类型"; Cell = "F24"; Old = 2989; New = 2991 },
    @{ Sheet = "全部类型"; Cell = "F25"; Old = 2637; New = 2638 },
    @{ Sheet = "全部类型"; Cell = "F31"; Old = 5363; New = 5366 },
    @{ Sheet = "全部类型"; Cell = "F34"; Old = 618;  New = 619 },
    @{ Sheet = "全部类型"; Cell = "F35"; Old = 618;  New = 619 },
    @{ Sheet = "全部类型"; Cell = "F38"; Old = 339;  New = 340 }
)

foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $cell = $ws.Range($change.Cell)
    $cell.Value = $change.New
}
